$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "244.80"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "21.75"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05770"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.413"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.299"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8169"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.043"
$ws.Range("E9").Value = "8FTXTokenFTT"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1427"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07312"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03124"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03136"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.134"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09448"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001596"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04829"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0005798"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006189"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.004138"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0009983"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0001498"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.731"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.160"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3264"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1292"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0003994"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03856"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006669"
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1072"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002740"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.006550"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005594"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.3895"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002097"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.01009"
